# Fruta / hortaliza, semanal
# Insert a new weekly record at row 8 (shifting existing rows 8-15 down to 9-16)
# for "Hortaliza, Mapocho Venta Directa de Santiago - Ají".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8, pushing rows 8-15 down to 9-16.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record's data.
$ws.Cells.Item(8, 1).Value = 12
$ws.Cells.Item(8, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(8, 3).Value = "Metropolitana"
$ws.Cells.Item(8, 4).Value = 44544
$ws.Cells.Item(8, 5).Value = 13
$ws.Cells.Item(8, 6).Value = 100112021
$ws.Cells.Item(8, 7).Value = "Ají"
$ws.Cells.Item(8, 8).Value = "Inferno"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 12
$ws.Cells.Item(8, 11).Value = 35000
$ws.Cells.Item(8, 12).Value = 35000
$ws.Cells.Item(8, 13).Value = 35000
$ws.Cells.Item(8, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 1400
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"
